$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0431000836353499
$ws.Range("C2").Value = 0.000446055199330917
$ws.Range("D2").Value = 0.000111513799832729
$ws.Range("E2").Value = 0.995539448006691
$ws.Range("F2").Value = 0.999498187900753
$ws.Range("G2").Value = 0.999609701700585
$ws.Range("H2").Value = 0.982938388625592
$ws.Range("I2").Value = 0.00100362419849456
$ws.Range("J2").Value = 0.998940618901589
$ws.Range("K2").Value = 0.13543350989685
$ws.Range("L2").Value = 0.998327293002509
$ws.Range("M2").Value = 0.000334541399498188
$ws.Range("N2").Value = 0.0218009478672986
$ws.Range("O2").Value = 0.999721215500418
$ws.Range("Q2").Value = 0.00161695009757457
$ws.Range("R2").Value = 0.999776972400335
$ws.Range("S2").Value = 0.99938667410092
$ws.Range("T2").Value = 0.0000557568999163647
$ws.Range("U2").Value = 0.000278784499581823
$ws.Range("V2").Value = 0.999052132701422
$ws.Range("W2").Value = 0.0248675773626986
$ws.Range("X2").Value = 0.00200724839698913
$ws.Range("B3").Value = 0.000446055199330917
$ws.Range("C3").Value = 0.0173961527739058
$ws.Range("D3").Value = 0.000446055199330917
$ws.Range("E3").Value = 0.0000557568999163647
$ws.Range("K3").Value = 0.0105938109841093
$ws.Range("L3").Value = 0.000111513799832729
$ws.Range("M3").Value = 0.999219403401171
$ws.Range("N3").Value = 0.904098132143853
$ws.Range("P3").Value = 0.000223027599665459
$ws.Range("Q3").Value = 0.948201839977697
$ws.Range("S3").Value = 0.0000557568999163647
$ws.Range("T3").Value = 0.00596598829105102
$ws.Range("U3").Value = 0.000111513799832729
$ws.Range("W3").Value = 0.000390298299414553
$ws.Range("B4").Value = 0.955394480066908
$ws.Range("C4").Value = 0.000669082798996376
$ws.Range("D4").Value = 0.0000557568999163647
$ws.Range("E4").Value = 0.00312238639531642
$ws.Range("F4").Value = 0.000501812099247282
$ws.Range("G4").Value = 0.000278784499581823
$ws.Range("H4").Value = 0.0156676888764985
$ws.Range("I4").Value = 0.998884862001673
$ws.Range("J4").Value = 0.000892110398661834
$ws.Range("K4").Value = 0.827822693058266
$ws.Range("L4").Value = 0.00150543629774185
$ws.Range("M4").Value = 0.000223027599665459
$ws.Range("N4").Value = 0.000167270699749094
$ws.Range("O4").Value = 0.000167270699749094
$ws.Range("P4").Value = 0.000278784499581823
$ws.Range("Q4").Value = 0.0000557568999163647
$ws.Range("R4").Value = 0.000167270699749094
$ws.Range("S4").Value = 0.000334541399498188
$ws.Range("T4").Value = 0.000278784499581823
$ws.Range("U4").Value = 0.999553944800669
$ws.Range("V4").Value = 0.000446055199330917
$ws.Range("W4").Value = 0.964761639252858
$ws.Range("X4").Value = 0.997936994703094
$ws.Range("B5").Value = 0.00083635349874547
$ws.Range("C5").Value = 0.980206300529691
$ws.Range("D5").Value = 0.99938667410092
$ws.Range("E5").Value = 0.00083635349874547
$ws.Range("G5").Value = 0.000111513799832729
$ws.Range("H5").Value = 0.0000557568999163647
$ws.Range("K5").Value = 0.0214664064678004
$ws.Range("L5").Value = 0.0000557568999163647
$ws.Range("M5").Value = 0.000223027599665459
$ws.Range("N5").Value = 0.0722051853916922
$ws.Range("O5").Value = 0.0000557568999163647
$ws.Range("P5").Value = 0.999442431000836
$ws.Range("Q5").Value = 0.04800669082799
$ws.Range("R5").Value = 0.0000557568999163647
$ws.Range("S5").Value = 0.000223027599665459
$ws.Range("T5").Value = 0.993364928909953
$ws.Range("U5").Value = 0.0000557568999163647
$ws.Range("V5").Value = 0.000390298299414553
$ws.Range("W5").Value = 0.0091441315862838
$ws.Range("X5").Value = 0.0000557568999163647
